{"js": "// Update the division-problem worksheet table: replace each problem's\n// text (\"N\u00f7M=\") with the new value, cell-by-cell, keeping run formatting\n// (rFonts/sz) untouched since we only change the cell's text value.\nconst replacements = [\n  [0, 0, \"52\u00f76=\", \"31\u00f73=\"],\n  [0, 1, \"23\u00f72=\", \"50\u00f78=\"],\n  [0, 2, \"15\u00f78=\", \"13\u00f72=\"],\n  [0, 3, \"43\u00f78=\", \"21\u00f72=\"],\n  [0, 4, \"79\u00f75=\", \"18\u00f78=\"],\n  [4, 0, \"50\u00f76=\", \"27\u00f78=\"],\n  [4, 1, \"35\u00f75=\", \"89\u00f77=\"],\n  [4, 2, \"74\u00f77=\", \"21\u00f76=\"],\n  [4, 3, \"24\u00f77=\", \"49\u00f75=\"],\n  [4, 4, \"53\u00f79=\", \"69\u00f75=\"],\n  [8, 0, \"58\u00f75=\", \"62\u00f76=\"],\n  [8, 1, \"69\u00f74=\", \"26\u00f78=\"],\n  [8, 2, \"18\u00f78=\", \"99\u00f74=\"],\n  [8, 3, \"81\u00f73=\", \"83\u00f76=\"],\n  [8, 4, \"69\u00f76=\", \"72\u00f76=\"],\n  [12, 0, \"27\u00f76=\", \"53\u00f77=\"],\n  [12, 1, \"40\u00f73=\", \"89\u00f72=\"],\n  [12, 2, \"22\u00f72=\", \"13\u00f72=\"],\n  [12, 3, \"35\u00f78=\", \"67\u00f79=\"],\n  [12, 4, \"55\u00f72=\", \"44\u00f73=\"],\n  [16, 0, \"43\u00f74=\", \"34\u00f72=\"],\n  [16, 1, \"78\u00f76=\", \"82\u00f77=\"],\n  [16, 2, \"78\u00f79=\", \"18\u00f78=\"],\n  [16, 3, \"71\u00f76=\", \"22\u00f72=\"],\n  [16, 4, \"41\u00f77=\", \"22\u00f74=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Load current values for every cell we intend to touch, so we can\n// confirm we are editing the expected problem before writing the new one.\nconst cells = replacements.map(([row, col]) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , oldText, newText] = replacements[i];\n  const cell = cells[i];\n  if (cell.value === oldText) {\n    cell.value = newText;\n  } else {\n    // Fallback: search-and-replace within this specific cell's body in\n    // case the cached value differs (defensive, should not trigger).\n    const searchResults = cell.body.search(oldText, { matchCase: true });\n    searchResults.load(\"items\");\n    await context.sync();\n    searchResults.items.forEach((range) => {\n      range.insertText(newText, Word.InsertLocation.replace);\n    });\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: replace each problem's\n# text (\"N\u00f7M=\") with the new value, cell-by-cell, keeping run formatting\n# (rFonts/sz) untouched since we only replace the text inside the cell's\n# existing range/run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Row = 1; Col = 1; Old = \"52\u00f76=\"; New = \"31\u00f73=\" },\n    @{ Row = 1; Col = 2; Old = \"23\u00f72=\"; New = \"50\u00f78=\" },\n    @{ Row = 1; Col = 3; Old = \"15\u00f78=\"; New = \"13\u00f72=\" },\n    @{ Row = 1; Col = 4; Old = \"43\u00f78=\"; New = \"21\u00f72=\" },\n    @{ Row = 1; Col = 5; Old = \"79\u00f75=\"; New = \"18\u00f78=\" },\n    @{ Row = 5; Col = 1; Old = \"50\u00f76=\"; New = \"27\u00f78=\" },\n    @{ Row = 5; Col = 2; Old = \"35\u00f75=\"; New = \"89\u00f77=\" },\n    @{ Row = 5; Col = 3; Old = \"74\u00f77=\"; New = \"21\u00f76=\" },\n    @{ Row = 5; Col = 4; Old = \"24\u00f77=\"; New = \"49\u00f75=\" },\n    @{ Row = 5; Col = 5; Old = \"53\u00f79=\"; New = \"69\u00f75=\" },\n    @{ Row = 9; Col = 1; Old = \"58\u00f75=\"; New = \"62\u00f76=\" },\n    @{ Row = 9; Col = 2; Old = \"69\u00f74=\"; New = \"26\u00f78=\" },\n    @{ Row = 9; Col = 3; Old = \"18\u00f78=\"; New = \"99\u00f74=\" },\n    @{ Row = 9; Col = 4; Old = \"81\u00f73=\"; New = \"83\u00f76=\" },\n    @{ Row = 9; Col = 5; Old = \"69\u00f76=\"; New = \"72\u00f76=\" },\n    @{ Row = 13; Col = 1; Old = \"27\u00f76=\"; New = \"53\u00f77=\" },\n    @{ Row = 13; Col = 2; Old = \"40\u00f73=\"; New = \"89\u00f72=\" },\n    @{ Row = 13; Col = 3; Old = \"22\u00f72=\"; New = \"13\u00f72=\" },\n    @{ Row = 13; Col = 4; Old = \"35\u00f78=\"; New = \"67\u00f79=\" },\n    @{ Row = 13; Col = 5; Old = \"55\u00f72=\"; New = \"44\u00f73=\" },\n    @{ Row = 17; Col = 1; Old = \"43\u00f74=\"; New = \"34\u00f72=\" },\n    @{ Row = 17; Col = 2; Old = \"78\u00f76=\"; New = \"82\u00f77=\" },\n    @{ Row = 17; Col = 3; Old = \"78\u00f79=\"; New = \"18\u00f78=\" },\n    @{ Row = 17; Col = 4; Old = \"71\u00f76=\"; New = \"22\u00f72=\" },\n    @{ Row = 17; Col = 5; Old = \"41\u00f77=\"; New = \"22\u00f74=\" }\n)\n\n$t = $d.Tables.Item(1)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $rng = $cell.Range\n    $rng.MoveEnd(1, -1) | Out-Null   # drop trailing cell-mark (wdCharacter = 1)\n\n    if ($rng.Text -eq $r.Old) {\n        $rng.Text = $r.New\n    } else {\n        # Fallback: scoped find/replace within this cell only.\n        $findRng = $cell.Range\n        $findRng.Find.ClearFormatting()\n        $findRng.Find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n    }\n}\n"}
